# Update the "TestData" worksheet:
#  1) Bump the date in C2 from 2025-06-18 (45826) to 2025-06-28 (45836).
#  2) Move the active cell/selection from B2 to G15 (last saved selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 holds a date (format yyyy-mm-dd) stored as the serial day number.
$ws.Range("C2").Value = 45836

# Update the worksheet's saved selection to G15.
$ws.Range("G15").Select()
